$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.23687527196074
$ws.Range("C2").Value = 10.97483311988117
$ws.Range("D2").Value = 4.625982044684009
$ws.Range("F2").Value = 24.49587423497817
$ws.Range("G2").Value = 3.62426158369357
$ws.Range("L2").Value = 10.67733898653682
$ws.Range("N2").Value = 17.23471492076004
$ws.Range("O2").Value = 21.76283323354574

$ws.Range("B3").Value = 15.71121012333023
$ws.Range("C3").Value = 10.85214211637441
$ws.Range("D3").Value = 4.599635227530389
$ws.Range("F3").Value = 24.44060631502906
$ws.Range("G3").Value = 3.626510406709676
$ws.Range("L3").Value = 10.6503850238575
$ws.Range("N3").Value = 17.29137618343138
$ws.Range("O3").Value = 21.78344730174181

$ws.Range("B4").Value = 15.38194902442727
$ws.Range("C4").Value = 10.77568604964794
$ws.Range("D4").Value = 4.583199294595151
$ws.Range("F4").Value = 24.41452419494837
$ws.Range("G4").Value = 3.627964938196083
$ws.Range("L4").Value = 10.63602602037264
$ws.Range("N4").Value = 17.32804434418824
$ws.Range("O4").Value = 21.80225419771837

$ws.Range("B5").Value = 15.24634991432424
$ws.Range("C5").Value = 10.74426732844416
$ws.Range("D5").Value = 4.576438691249759
$ws.Range("F5").Value = 24.40587674953144
$ws.Range("G5").Value = 3.628576275378022
$ws.Range("L5").Value = 10.63072963986819
$ws.Range("N5").Value = 17.3434602235942
$ws.Range("O5").Value = 21.81146082610087

$ws.Range("B6").Value = 15.22375439311891
$ws.Range("C6").Value = 10.73903502981833
$ws.Range("D6").Value = 4.575312369946895
$ws.Range("F6").Value = 24.40456065568589
$ws.Range("G6").Value = 3.62867891283972
$ws.Range("L6").Value = 10.62988381254986
$ws.Range("N6").Value = 17.3460486374641
$ws.Range("O6").Value = 21.81308262721674

$ws.Range("B7").Value = 15.38012575468908
$ws.Range("C7").Value = 10.77526336075663
$ws.Range("D7").Value = 4.58310837013037
$ws.Range("F7").Value = 24.41439954362041
$ws.Range("G7").Value = 3.627973107497864
$ws.Range("L7").Value = 10.63595233917926
$ws.Range("N7").Value = 17.32825033028614
$ws.Range("O7").Value = 21.80237212093565

$ws.Range("B8").Value = 16.05709714921334
$ws.Range("C8").Value = 10.93277443372141
$ws.Range("D8").Value = 4.616952751421771
$ws.Range("F8").Value = 24.47519274913984
$ws.Range("G8").Value = 3.62502170824209
$ws.Range("L8").Value = 10.66759318270069
$ws.Range("N8").Value = 17.25386261614936
$ws.Range("O8").Value = 21.76866257153377

$ws.Range("B9").Value = 17.32460043766047
$ws.Range("C9").Value = 11.23187064095501
$ws.Range("D9").Value = 4.681169084101947
$ws.Range("F9").Value = 24.65631283181272
$ws.Range("G9").Value = 3.619816421225496
$ws.Range("L9").Value = 10.74682013720139
$ws.Range("N9").Value = 17.12283949794999
$ws.Range("O9").Value = 21.75148556035845

$ws.Range("B10").Value = 18.20911698939092
$ws.Range("C10").Value = 11.44451778795506
$ws.Range("D10").Value = 4.726907187550234
$ws.Range("F10").Value = 24.82642275188659
$ws.Range("G10").Value = 3.616343315108737
$ws.Range("L10").Value = 10.81519562968939
$ws.Range("N10").Value = 17.03556146753696
$ws.Range("O10").Value = 21.76881966905951

$ws.Range("B11").Value = 18.59955789393729
$ws.Range("C11").Value = 11.53948323005029
$ws.Range("D11").Value = 4.747374341937801
$ws.Range("F11").Value = 24.91165515013332
$ws.Range("G11").Value = 3.614838758743937
$ws.Range("L11").Value = 10.8484336928371
$ws.Range("N11").Value = 16.99779303224632
$ws.Range("O11").Value = 21.78321764058817

$ws.Range("B12").Value = 18.74556309148964
$ws.Range("C12").Value = 11.57517100744043
$ws.Range("D12").Value = 4.755073603953459
$ws.Range("F12").Value = 24.94503881192478
$ws.Range("G12").Value = 3.614279799678875
$ws.Range("L12").Value = 10.86131975963754
$ws.Range("N12").Value = 16.98376829724988
$ws.Range("O12").Value = 21.78960526726143

$ws.Range("B13").Value = 18.71420230044373
$ws.Range("C13").Value = 11.5674974729039
$ws.Range("D13").Value = 4.753417745596193
$ws.Range("F13").Value = 24.93780011533229
$ws.Range("G13").Value = 3.614399702780011
$ws.Range("L13").Value = 10.85853131359012
$ws.Range("N13").Value = 16.98677645310612
$ws.Range("O13").Value = 21.78818800187236

$ws.Range("B14").Value = 18.61160752868518
$ws.Range("C14").Value = 11.54242487574738
$ws.Range("D14").Value = 4.748008804803751
$ws.Range("F14").Value = 24.91437956375841
$ws.Range("G14").Value = 3.614792557038012
$ws.Range("L14").Value = 10.84948787755806
$ws.Range("N14").Value = 16.99663365519222
$ws.Range("O14").Value = 21.78372441681467

$ws.Range("B15").Value = 18.5485211969485
$ws.Range("C15").Value = 11.52703099845916
$ws.Range("D15").Value = 4.744688932853547
$ws.Range("F15").Value = 24.90017744581385
$ws.Range("G15").Value = 3.615034594344832
$ws.Range("L15").Value = 10.84398729594986
$ws.Range("N15").Value = 17.00270757017515
$ws.Range("O15").Value = 21.78111211345996

$ws.Range("B16").Value = 18.1833491759099
$ws.Range("C16").Value = 11.43827432312125
$ws.Range("D16").Value = 4.725562576128741
$ws.Range("F16").Value = 24.82100883494294
$ws.Range("G16").Value = 3.61644315344273
$ws.Range("L16").Value = 10.81306576683751
$ws.Range("N16").Value = 17.03806858160247
$ws.Range("O16").Value = 21.76800970850168

$ws.Range("B17").Value = 17.95617458613125
$ws.Range("C17").Value = 11.38335807711604
$ws.Range("D17").Value = 4.713740542658536
$ws.Range("F17").Value = 24.77443768835558
$ws.Range("G17").Value = 3.617326525176424
$ws.Range("L17").Value = 10.79463786131924
$ws.Range("N17").Value = 17.0602563350648
$ws.Range("O17").Value = 21.7616393886777

$ws.Range("B18").Value = 17.82439368913069
$ws.Range("C18").Value = 11.35160632895149
$ws.Range("D18").Value = 4.706909024377916
$ws.Range("F18").Value = 24.74839131386407
$ws.Range("G18").Value = 3.617841715356809
$ws.Range("L18").Value = 10.78423991815762
$ws.Range("N18").Value = 17.07320030471206
$ws.Range("O18").Value = 21.75858847281352

$ws.Range("B19").Value = 17.77958748269537
$ws.Range("C19").Value = 11.34082792459674
$ws.Range("D19").Value = 4.704590602199755
$ws.Range("F19").Value = 24.73970016707292
$ws.Range("G19").Value = 3.618017370704131
$ws.Range("L19").Value = 10.78075414156664
$ws.Range("N19").Value = 17.07761422355248
$ws.Range("O19").Value = 21.75766080631066

$ws.Range("B20").Value = 17.9804742135387
$ws.Range("C20").Value = 11.38922126219121
$ws.Range("D20").Value = 4.715002326713918
$ws.Range("F20").Value = 24.77931880941187
$ws.Range("G20").Value = 3.617231754520263
$ws.Range("L20").Value = 10.79657876083921
$ws.Range("N20").Value = 17.05787556422845
$ws.Range("O20").Value = 21.76225406872376

$ws.Range("B21").Value = 18.64179317726596
$ws.Range("C21").Value = 11.54979688019764
$ws.Range("D21").Value = 4.74959895060776
$ws.Range("F21").Value = 24.92122885152143
$ws.Range("G21").Value = 3.614676873989819
$ws.Range("L21").Value = 10.85213608500177
$ws.Range("N21").Value = 16.99373083553415
$ws.Range("O21").Value = 21.78501010621012

$ws.Range("B22").Value = 19.06319122488627
$ws.Range("C22").Value = 11.65313924934476
$ws.Range("D22").Value = 4.771910280583334
$ws.Range("F22").Value = 25.02042250037508
$ws.Range("G22").Value = 3.613069942585156
$ws.Range("L22").Value = 10.89018822783218
$ws.Range("N22").Value = 16.9534247812041
$ws.Range("O22").Value = 21.80533367810715

$ws.Range("B23").Value = 18.83931190913022
$ws.Range("C23").Value = 11.59813635482749
$ws.Range("D23").Value = 4.760030503825228
$ws.Range("F23").Value = 24.96689841962952
$ws.Range("G23").Value = 3.613921861282411
$ws.Range("L23").Value = 10.8697221926542
$ws.Range("N23").Value = 16.97478927537881
$ws.Range("O23").Value = 21.79398843355983

$ws.Range("B24").Value = 17.96949200236309
$ws.Range("C24").Value = 11.38657107439299
$ws.Range("D24").Value = 4.714431982572047
$ws.Range("F24").Value = 24.7771097853434
$ws.Range("G24").Value = 3.617274577503058
$ws.Range("L24").Value = 10.7957006675308
$ws.Range("N24").Value = 17.0589513252288
$ws.Range("O24").Value = 21.76197426708743

$ws.Range("B25").Value = 16.98924504394755
$ws.Range("C25").Value = 11.15212332692213
$ws.Range("D25").Value = 4.664040689753968
$ws.Range("F25").Value = 24.60075028677333
$ws.Range("G25").Value = 3.621162636745784
$ws.Range("L25").Value = 10.72357748912665
$ws.Range("N25").Value = 17.15670203904389
$ws.Range("O25").Value = 21.75087766831054
